$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and the worksheet's own Name property) from 2022-02-20 -> 2022-02-21
$ws.Name = "Through 2022-02-21"

# Update the "February (through 02-20)" label to "February (through 02-21)"
$ws.Range("A3").Value = "February (through 02-21)"

# Update row 3 (February) values
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 41
$ws.Range("F3").Value = 22
$ws.Range("G3").Value = 56
$ws.Range("H3").Value = 96
$ws.Range("I3").Value = 106

# Update row 4 (Total) values
$ws.Range("B4").Value = 35
$ws.Range("C4").Value = 77
$ws.Range("D4").Value = 117
$ws.Range("E4").Value = 127
$ws.Range("F4").Value = 71
$ws.Range("G4").Value = 130
$ws.Range("H4").Value = 313
$ws.Range("I4").Value = 265

$wb.Save()
